$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update "想去人数" (F column) for three rows
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 230
$wsExhibit.Range("F4").Value = 847
$wsExhibit.Range("F6").Value = 33

# Sheet "全部类型" (All types) - same three events shifted down by one row
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 230
$wsAll.Range("F5").Value = 847
$wsAll.Range("F7").Value = 33
